# Updated cryptos list - applies the price/volume/coin changes described in the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.830.12"
$ws.Range("E2").Value = "  -1.57%  "
$ws.Range("D3").Value = "2.659.96"
$ws.Range("E3").Value = "  +0.55%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'598.24"
$ws.Range("E5").Value = "  -1.41%  "
$ws.Range("D6").Value = "'174.11"
$ws.Range("E6").Value = "  -2.95%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").Value = "'0.523"
$ws.Range("E8").Value = "  -0.90%  "
$ws.Range("D9").Value = "2.659.45"
$ws.Range("E9").Value = "  +0.56%  "
$ws.Range("E10").Value = "  -2.74%  "
$ws.Range("E11").Value = "  +2.34%  "
$ws.Range("E12").Value = "  +0.35%  "
$ws.Range("E13").Value = "  -1.78%  "
$ws.Range("D14").Value = "3.148.67"
$ws.Range("E14").Value = "  +0.33%  "
$ws.Range("E15").Value = "  -3.21%  "
$ws.Range("D16").Value = "71.842.67"
$ws.Range("E16").Value = "  -1.57%  "
$ws.Range("D17").Value = "'26.18"
$ws.Range("E17").Value = "  -2.36%  "
$ws.Range("D18").Value = "2.664.48"
$ws.Range("E18").Value = "  -1.24%  "
$ws.Range("E19").Value = "  +5.43%  "
$ws.Range("D20").Value = "'8.18"
$ws.Range("E20").Value = "  +1.70%  "
$ws.Range("D21").Value = "'370.95"
$ws.Range("D22").Value = "'4.15"
$ws.Range("E22").Value = "  -1.27%  "
$ws.Range("E23").Value = "  +0.65%  "
$ws.Range("D24").Value = "'71.89"
$ws.Range("E24").Value = "  -1.87%  "
$ws.Range("E25").Value = "  +0.00%  "
$ws.Range("D26").Value = "'4.32"
$ws.Range("E26").Value = "  -1.74%  "
$ws.Range("D27").Value = "'9.73"
$ws.Range("E27").Value = "  -1.70%  "
$ws.Range("D28").Value = "2.797.56"
$ws.Range("E28").Value = "  +0.29%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  -0.20%  "
$ws.Range("D30").Value = "0.0₃0964"
$ws.Range("E30").Value = "  -0.30%  "
$ws.Range("D31").Value = "'8.04"
$ws.Range("E31").Value = "  -0.30%  "
$ws.Range("D32").Value = "'499.64"
$ws.Range("E32").Value = "  -6.48%  "
$ws.Range("E33").Value = "  -3.38%  "
$ws.Range("E34").Value = "  -0.84%  "
$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("D36").Value = "'163.01"
$ws.Range("E36").Value = "  -0.32%  "
$ws.Range("D37").Value = "'19.46"
$ws.Range("E37").Value = "  +0.50%  "
$ws.Range("D38").Value = "'19.06"
$ws.Range("E38").Value = "  -0.28%  "
$ws.Range("E39").Value = "  -2.81%  "
$ws.Range("E40").Value = "  -2.71%  "
$ws.Range("E41").Value = "  -4.40%  "
$ws.Range("E42").Value = "  -0.05%  "
$ws.Range("D43").Value = "'4.98"
$ws.Range("E43").Value = "  -2.56%  "
$ws.Range("E44").Value = "  -3.10%  "
$ws.Range("E45").Value = "  -0.62%  "
$ws.Range("B46").Value = "OKB"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D46").Value = "'39.47"
$ws.Range("E46").Value = "  -0.56%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "'155.96"
$ws.Range("E47").Value = "  +2.65%  "
$ws.Range("D48").Value = "'0.558"
$ws.Range("E48").Value = "  +2.44%  "
$ws.Range("D49").Value = "'3.72"
$ws.Range("E49").Value = "  +0.60%  "
$ws.Range("E50").Value = "  +1.18%  "
$ws.Range("D51").Value = "'0.0752"
$ws.Range("E51").Value = "  -1.82%  "
